# Auto-generated edit script: "Generate Report for Handoff"
# Adds a new localization entry (2861cb3a-...md) to the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$newMdFile   = "2861cb3a-6b60-44da-ac15-ae357a2c6da3ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdDisplay = "e2e\2861cb3a-6b60-44da-ac15-ae357a2c6da3ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$statusReady = "Ready for handoff"
$date51      = "2016-08-16 16:24:51"
$newZhCnXlf  = "2861cb3a-6b60-44da-ac15-ae357a2c6da3oooooooooooooooooooooooooooooooooooooooo.8e4df41575d7c4bd2937950fdf3a3e135818b3c3.zh-cn.xlf"
$date46      = "2016-08-16 16:24:46"
$newDeDeXlf  = "2861cb3a-6b60-44da-ac15-ae357a2c6da3oooooooooooooooooooooooooooooooooooooooo.8e4df41575d7c4bd2937950fdf3a3e135818b3c3.de-de.xlf"
$newUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd40cad7e47d133c0135454a8bf125322af94ea3/e2e/2861cb3a-6b60-44da-ac15-ae357a2c6da3ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

$hyperlinkColor = 15570276  # BGR int for RGB 6495ED (matches existing HyperLink style)

function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Set-DateLook($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview": add row 3
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = $newMdFile
$wsOv.Range("B3").Value = $newMdDisplay
Set-HyperlinkLook($wsOv.Range("B3"))
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = "'"
$wsOv.Range("E3").Value = $statusReady
$wsOv.Range("F3").Value = $statusReady
$wsOv.Range("G3").Value = $date51
Set-DateLook($wsOv.Range("G3"))

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $newUrl, "", "", $newMdDisplay) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newMdFile
Set-HyperlinkLook($wsZh.Range("A3"))
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $newZhCnXlf
$wsZh.Range("H3").Value = $date46
Set-DateLook($wsZh.Range("H3"))
$wsZh.Range("I3").Value = "'"
$wsZh.Range("J3").Value = "'"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
Set-DateLook($wsZh.Range("K3"))
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newMdFile) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": add row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newMdFile
Set-HyperlinkLook($wsDe.Range("A3"))
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $newDeDeXlf
$wsDe.Range("H3").Value = $date51
Set-DateLook($wsDe.Range("H3"))
$wsDe.Range("I3").Value = "'"
$wsDe.Range("J3").Value = "'"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
Set-DateLook($wsDe.Range("K3"))
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newMdFile) | Out-Null

Write-Output "Report generated for handoff"
